{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that reads \"\u2022 Expert methodology validated at highest judicial level\"\n// and insert the two new achievement paragraphs immediately after it.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"\u2022 Expert methodology validated at highest judicial level\") {\n    anchor = p;\n    break;\n  }\n}\n\nif (anchor) {\n  // First new paragraph: plain text only.\n  const p1 = anchor.insertParagraph(\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n    \"After\"\n  );\n\n  // Second new paragraph: \"\u2022 \" + bold/colored \"178%\" + \" accuracy improvement in racial classification algorithms\"\n  const p2 = p1.insertParagraph(\"\u2022 \", \"After\");\n  const r2 = p2.insertText(\"178%\", \"End\");\n  r2.font.bold = true;\n  r2.font.color = \"#2C3E50\";\n  p2.insertText(\" accuracy improvement in racial classification algorithms\", \"End\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph \"\u2022 Expert methodology validated at highest judicial level\"\n$targetIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"\u2022 Expert methodology validated at highest judicial level\") {\n        $targetIndex = $i\n        break\n    }\n    $i = $i + 1\n}\n\nif ($targetIndex -gt 0) {\n    $target = $d.Paragraphs($targetIndex)\n\n    # Insert a new empty paragraph right after the target paragraph, then fill it in.\n    $target.Range.InsertParagraphAfter()\n    $p1 = $d.Paragraphs($targetIndex + 1)\n    $p1.Range.Text = \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n\n    # Insert a second new empty paragraph right after the first new one.\n    $p1.Range.InsertParagraphAfter()\n    $p2 = $d.Paragraphs($targetIndex + 2)\n    $p2.Range.Text = \"\u2022 \"\n\n    # Append the bold/colored \"178%\" run right after \"\u2022 \".\n    $boldRange = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)\n    $boldRange.InsertAfter(\"178%\")\n    $boldRange.Font.Bold = 1\n    $boldRange.Font.Color = 5258796\n\n    # Append the trailing plain-text run.\n    $tailRange = $d.Range($boldRange.End, $boldRange.End)\n    $tailRange.InsertAfter(\" accuracy improvement in racial classification algorithms\")\n}\n"}
